# 8th day of DSA Series
# Update workbook: mark several "Searching & Sorting" rows as solved/unsolved,
# add a new "Searching & Sorting *" category marker, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Section separator rows: append " *" to the topic name once a sub-item below
# it has been revisited.
# ---------------------------------------------------------------------------
$ws.Cells.Item(61, 1).Value = "String *"
$ws.Cells.Item(80, 1).Value = "String *"
$ws.Cells.Item(106, 1).Value = "Searching & Sorting *"

# ---------------------------------------------------------------------------
# Rows 115-130: mark each question as done ("yes", green link) or not done
# ("no", red link + topic marked with "Searching & Sorting *"). Row 116 is
# marked with a distinct (new) highlight color.
# ---------------------------------------------------------------------------

function Set-Done($row) {
    $c = $ws.Cells.Item($row, 2)
    $c.Font.Underline = $true
    $c.Font.Color = 5287936   # RGB(0,176,80) -- green
    $ws.Cells.Item($row, 3).Value = "yes"
}

function Set-NotDone($row) {
    $ws.Cells.Item($row, 1).Value = "Searching & Sorting *"
    $c = $ws.Cells.Item($row, 2)
    $c.Font.Underline = $true
    $c.Font.Color = 255       # RGB(255,0,0) -- red
    $ws.Cells.Item($row, 3).Value = "no"
}

function Set-DoneHighlight($row) {
    $c = $ws.Cells.Item($row, 2)
    $c.Font.Underline = $true
    $c.Font.ThemeColor = 6    # -> theme="5" (accent2) in the saved style
    $ws.Cells.Item($row, 3).Value = "yes"
}

Set-Done 115
Set-DoneHighlight 116
Set-Done 117
Set-NotDone 118
Set-Done 119
Set-NotDone 120
Set-Done 121
Set-Done 122
Set-Done 123
Set-Done 124
Set-Done 125
Set-Done 126
Set-NotDone 127
Set-Done 128
Set-Done 129
Set-Done 130

# ---------------------------------------------------------------------------
# Move the view / selection to where we were working.
# ---------------------------------------------------------------------------
$ws.Range("A118").Select()
$excel.ActiveWindow.ScrollRow = 118
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C127").Select()
